$wb = $excel.ActiveWorkbook

# --- Variants sheet: insert a new "Description" column before the old
# "Relation Type(s)" column (old G -> new H), and fill in two description
# values for the first two data rows. ---
$ws = $wb.Worksheets.Item("Variants")

$ws.Columns("G:G").Insert()

$ws.Range("G1").Value = "Description"
$ws.Range("G2").Value = "description1"
$ws.Range("G3").Value = "description2"

# Bold header font keeps the same font id, so this merely adds the
# missing <family val="2"/> to the existing bold font definition.
$ws.Range("A1:H1").Font.Name = "Calibri"

# Page setup metadata added by the save.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection moved.
$ws.Range("I10").Select() | Out-Null

# --- Categories sheet: COUNTIF ranges shift from G:P to H:Q to track the
# inserted column on Variants. ---
$cs = $wb.Worksheets.Item("Categories")

$cs.Range("C2").Formula = '=COUNTIF(Variants!H:Q, "category1")'
$cs.Range("D2").Formula = '=COUNTIF(Variants!H:Q, "category1")'

$cs.Range("C3").Formula = '=COUNTIF(Variants!H:Q, "category2")'
$cs.Range("D3").Formula = '=COUNTIF(Variants!H:Q, "category2")'

$cs.Range("C4").Formula = '=COUNTIF(Variants!H:Q, "category3")'
$cs.Range("D4").Formula = '=COUNTIF(Variants!H:Q, "category3")'
